{"js": "// HTH Chapter 17 - fix \"awakened\" -> \"woke up\" (x2) and \"mistake one\" -> \"mistake once\"\nconst body = context.document.body;\n\nconst replacements = [\n  {\n    find: \"Gayoon awakened, weakly opening her eyes.\",\n    replace: \"Gayoon woke up, weakly opening her eyes.\"\n  },\n  {\n    find: \"\\\"You awakened...?\\\" - A voice behind her spoke.\",\n    replace: \"\\\"You woke up...?\\\" - A voice behind her spoke.\"\n  },\n  {\n    find: \"I made this mistake one and I don't want you to repeat it\",\n    replace: \"I made this mistake once and I don't want you to repeat it\"\n  }\n];\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${find}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# HTH Chapter 17 - fix \"awakened\" -> \"woke up\" (x2) and \"mistake one\" -> \"mistake once\"\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @('Gayoon awakened, weakly opening her eyes.', 'Gayoon woke up, weakly opening her eyes.'),\n    @('\"You awakened...?\" - A voice behind her spoke.', '\"You woke up...?\" - A voice behind her spoke.'),\n    @(\"I made this mistake one and I don't want you to repeat it\", \"I made this mistake once and I don't want you to repeat it\")\n)\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\n$d.Content.Text.Substring(0, 60)\n"}
